$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row: copy O1 format into P1:Q1, then set the two new header values ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- New columns P and Q hold zeros for every data row (2-25) ---
for ($r = 2; $r -le 25; $r++) {
  $ws.Cells.Item($r, 16).Value = 0   # column P
  $ws.Cells.Item($r, 17).Value = 0   # column Q
}

# --- Updated simulation results for columns B:K (2..11), rows 2-25 ---
$data = @(
  @(3.559570922419823, 1.322763020649631, 0.01812025240789694, 1.433382523723949, 0.7870434807789479, 0.7110333397667432, 0, 0.0157539186156006, 0.4252747715031262, 0.2552068450048992),
  @(3.094976390801605, 1.152411130016333, 0.01716004229405854, 1.227872290318089, 0.7036598521921462, 0.6314952471928166, 0, 0.01052940931996149, 0.3964458630915431, 0.2499594840145143),
  @(2.810397383663428, 1.049031273224699, 0.0165927334427316, 1.10446850180027, 0.6540834212456446, 0.5843015354369214, 0, 0.007820485930202903, 0.3796729111232935, 0.2474710889212659),
  @(2.69449533423699, 1.008455067574658, 0.01642172002073039, 1.054760344638098, 0.6329936309558377, 0.5640452183546074, 0, 0.006877645738713856, 0.3723358141636055, 0.2458695820557573),
  @(2.675196263437613, 1.003313727383045, 0.01646128554909154, 1.046526575455871, 0.6279861951028423, 0.5590062995076295, 0, 0.006793735659964995, 0.3702585280635162, 0.2446922092085408),
  @(2.808669164602918, 1.05281582574284, 0.01677590525476802, 1.103760171294837, 0.6496200408859352, 0.5793829823974903, 0, 0.007971119807720584, 0.3771905170936662, 0.2449294551226764),
  @(3.398965382442043, 1.269542480187511, 0.01803444382436226, 1.361820312999839, 0.7523643548725545, 0.6770496423210659, 0, 0.01402166165024532, 0.4119640928718411, 0.249897644139299),
  @(4.565207159732608, 1.698555288222053, 0.02035984581094752, 1.895573088683832, 0.9765557790446024, 0.8924798321767611, 0, 0.03024497693943218, 0.4929276059901753, 0.2708083691579191),
  @(5.431640710337831, 2.026699701013797, 0.0223214332310171, 2.313425999797417, 1.151988948629736, 1.061735550352424, 0, 0.04597216024347794, 0.5581447946724154, 0.29026973540185),
  @(5.828710696860298, 2.186635384131534, 0.02358907277619338, 2.511356488424198, 1.228487369889478, 1.134821551820778, 0, 0.05438770208917987, 0.5858138146548697, 0.2964970713215962),
  @(5.979848623111479, 2.243819626842992, 0.02390154104429598, 2.58775041377082, 1.26178962720941, 1.167281842275116, 0, 0.05763683939692044, 0.5987123112562074, 0.3012891776991253),
  @(5.947310475330653, 2.230683526781888, 0.02379835149955767, 2.571245110281581, 1.255337100347958, 1.161092710237682, 0, 0.05690415461790455, 0.596341485377053, 0.3006850317385243),
  @(5.841149328432152, 2.190991159680664, 0.02359954031611622, 2.517617732819289, 1.231531970223173, 1.137831654563144, 0, 0.05464046849792048, 0.5870474514299531, 0.2970726947121207),
  @(5.776125075151072, 2.168315199718108, 0.02354866663868904, 2.484932309792711, 1.215563125008657, 1.12203661729032, 0, 0.05333109522536006, 0.580568323563142, 0.2940299796300962),
  @(5.405132280994849, 2.029432525131313, 0.02282294718022371, 2.300473667600357, 1.134685969563392, 1.043273918805568, 0, 0.04585086898000057, 0.5493861179907356, 0.2826022045112424),
  @(5.17846635025694, 1.945051936881384, 0.02238644645592203, 2.189515356639859, 1.086140057801558, 0.9960848484995779, 0, 0.04155654236696815, 0.5308203927576187, 0.2760307474189645),
  @(5.048583759827466, 1.892828799951928, 0.02196123477039791, 2.126484971250974, 1.062311876077501, 0.9734709161204762, 0, 0.03902254091252111, 0.5224418628185106, 0.2745930252111535),
  @(5.004514855737909, 1.878250146643666, 0.02195426522393262, 2.105221206079449, 1.051426151246218, 0.9626793925862387, 0, 0.0383047764823008, 0.5180110671929583, 0.2724397322075731),
  @(5.202577534689794, 1.953693566231152, 0.02241864663460547, 2.201253792187757, 1.091562411417968, 1.001393226543073, 0, 0.04198997446944475, 0.5329418699809878, 0.276886318055233),
  @(5.872092091069192, 2.20709231327919, 0.02385761356082128, 2.533257239903079, 1.234311598197522, 1.139982478609852, 0, 0.05545245276405364, 0.5874026352450414, 0.2956677168937105),
  @(6.313431476763583, 2.370065536351603, 0.02457572439544364, 2.758538355958905, 1.336809399239826, 1.240584110122001, 0, 0.06511286364086732, 0.6280296226873361, 0.3126382603158433),
  @(6.077793306233843, 2.277401592166598, 0.02394859959053264, 2.637552521886064, 1.286782698585753, 1.192112060138157, 0, 0.0596876104866455, 0.6090004209187896, 0.3063967478822107),
  @(5.192072117481302, 1.941392245765883, 0.0220300692109241, 2.196064314732595, 1.097022212852764, 1.007786779737089, 0, 0.04148837377826009, 0.5364515066278273, 0.2811541889428995),
  @(4.248013901697732, 1.588934612195771, 0.02005918487907721, 1.747279953126437, 0.9062846321032794, 0.8236573949315158, 0, 0.025526115399364, 0.4656119911405341, 0.2597694914840503),
)

$startRow = 2
$startCol = 2  # column B
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowArr = $data[$i]
  for ($j = 0; $j -lt $rowArr.Count; $j++) {
    $c = $startCol + $j
    $ws.Cells.Item($r, $c).Value = $rowArr[$j]
  }
}
